$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two rows (19 and 20), shrinking the used range to A1:C18
$ws.Rows("19:20").Delete()

# Replace the transcript text in column C (rows 2-18) with the new content
$ws.Range("C2").Value = "めちゃキモいよ これさ 俺の 入ってるこれ 声の路地裏の声 入ってくらに"
$ws.Range("C3").Value = "頭大丈夫と思う バカゲー"
$ws.Range("C4").Value = "柴崎 ゴルフ"
$ws.Range("C5").Value = "難波まで言われたら絶対的な権力が欲しいかな でも なんか なんじゃないかな 俺は河川 河川とゲームをすることから 一緒に素晴らしい一番 現実味があり そうやけど 俺の名前 見たことは"
$ws.Range("C6").Value = "今まで見たことはあるからそんなあるおしっこ 俺 現れ 一緒に あれか"
$ws.Range("C7").Value = "ランボルギーニ 欲しいところにない車になりたい所があって 明日車になってから喋れへんとか苦労して"
$ws.Range("C8").Value = "僕は過去に戻れたらええかな 終わりや それもあるんやん 俺とかなんか人生とかなんかやり直したらええなぁとかと思うんやけどなんかこれから先とかあるんじゃないなんか ここに行きたいと言うか なんかそれもあるけどなんかどっかに挫折したところがあったとしたらなんかその挫折する前に戻れたらええなぁとかと思ってやり直したいなーとか白いな やらないことだらけで今すぐやり直したいことじゃないんや そうじゃないけど あれかいや でもこれからさ"
$ws.Range("C9").Value = "ツェナー 願い事がひとつ叶うなら やけど 一回一回だけってことで起こった権力があったら何でも手に入るやん このようなものは金もそう やしま コネとか関係の引力に群がる人間と すぐもう 跪かせられる 麻疹に絶対的な権力があったら投資もできるやん 自分のやりたいように 世界を動かすこともできる しにく 訳あり 件とかもすごい 晴れたりできる ロックフェラー"
$ws.Range("C10").Value = "ロックフェラー"
$ws.Range("C11").Value = "筋力あったらいける方もいるけど筋力で解決するのかな と思ったけどまあなんか 宇宙好き だし 俺は強かっただけで ただ ただ あの宇宙から地球を見てみたいのがあるよね 今から死ぬほど勉強しやんといてよろし 宇宙 なんて 元でどうぞ"
$ws.Range("C12").Value = "どこにそう どこにある 感じたがいきなりだとカリスマ性からそのせいか 夜になると思うんだったらよかったなあ なんかやったよ これだと一緒にそっちなくて有名な"
$ws.Range("C13").Value = "何してもわかれば許されるし 結局 妹 が なんかしなあかんとなるけど育ったお金だけが手に入るっていう 札がやったらさ一生労働しなくていいっちゅう メリットってあれはな 何だ それで 腕"
$ws.Range("C14").Value = "もしかしたら金が無限にまる だからなんかなんかその代わりじゃないなんか大嫌い 特に 4月とかになったらお金を可能なるって何が原因"
$ws.Range("C15").Value = "田舎 お金あっても 宇宙お金あったら 家 行けると行けるのも 起こりつつある いけんのかに何ができる"
$ws.Range("C16").Value = "過去に戻れるかな 人生の過去じゃなくてもっと過去に戻りたいのあるかもしれ 信長の時代とかそういうの見れるんだったらそっちの方がいいかも 自分の"
$ws.Range("C17").Value = "現在地の知識で石器時代とか言ったらさー 存在その子が神として崇拝されるレベルだけど多分先へ来るんですか もう1回だけやからもう戻るんやったらもう 石器時代の石器時代でもそろそろ死ぬまで生きれるかどうもこれ まとめ 結論 お金で全てが解決できるけどま未来に戻ることはできんから一概には金が全てとは言えへん"
$ws.Range("C18").Value = "ところ これからしか来ないんでしょ あと30秒ぐらい 15秒"
